$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing trend values for P5:R5 (2019-2021 columns)
$ws.Range("P5").Value = 4.4000000000000004
$ws.Range("Q5").Value = 2.9
$ws.Range("R5").Value = 3.2

# Add a new column S (year 2022) by copying the formatting of column R
# (which carries the correct style indexes for the header row and the
# data row) one column to the right, then filling in the new values.
$ws.Range("R4:R5").Copy() | Out-Null
$ws.Range("S4:S5").Insert(-4161) | Out-Null

$ws.Range("S4").Value = 2022
$ws.Range("S5").Value = 3.4

# Move the active selection the way the author's workbook shows it
$ws.Range("T4").Select() | Out-Null
